$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''55.630.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.61%  '
$ws.Range("D3").Value = '''2.484.98'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.40%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '''481.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.01%  '
$ws.Range("D6").Value = '''143.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +10.55%  '
$ws.Range("D7").Value = '''0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("D8").Value = '''0.506'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.21%  '
$ws.Range("D9").Value = '''2.506.31'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.62%  '
$ws.Range("D10").Value = '''5.63'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.26%  '
$ws.Range("D11").Value = '''0.0966'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.79%  '
$ws.Range("D12").Value = '''0.329'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.09%  '
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("D14").Value = '''2.915.63'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.57%  '
$ws.Range("D15").Value = '''55.822.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.45%  '
$ws.Range("D16").Value = '''20.87'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.74%  '
$ws.Range("D17").Value = '''0.0000135'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.08%  '
$ws.Range("D18").Value = '''2.508.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.98%  '
$ws.Range("D19").Value = '''4.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.93%  '
$ws.Range("E20").Value = '  +9.74%  '
$ws.Range("D21").Value = '''317.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.44%  '
$ws.Range("D22").Value = '''0.997'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").Value = '''5.76'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.56%  '
$ws.Range("D24").Value = '''58.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.08%  '
$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").Value = '''0.166'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.58%  '
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").Value = '''0.407'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.22%  '
$ws.Range("D27").Value = '''0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.51%  '
$ws.Range("D28").Value = '''2.611.31'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.66%  '
$ws.Range("D29").Value = '''7.46'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.93%  '
$ws.Range("D30").Value = '''0.0₃0776'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.27%  '
$ws.Range("E31").Value = '  +0.31%  '
$ws.Range("D32").Value = '''148.38'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.89%  '
$ws.Range("D33").Value = '''18.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.85%  '
$ws.Range("D34").Value = '''1.48'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.20%  '
$ws.Range("D35").Value = '''5.19'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.53%  '
$ws.Range("E36").Value = '  +9.19%  '
$ws.Range("D37").Value = '''3.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.15%  '
$ws.Range("D38").Value = '''0.855'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.37%  '
$ws.Range("D39").Value = '''34.10'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.10%  '
$ws.Range("D40").Value = '''3.53'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.34%  '
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '''0.996'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.32%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Value = '''0.609'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.89%  '
$ws.Range("D43").Value = '''0.0552'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.45%  '
$ws.Range("D44").Value = '''1.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.50%  '
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").Value = '''258.90'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +20.80%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").Value = '''4.68'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.56%  '
$ws.Range("D47").Value = '''10.19'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '''0.0900'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.27%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '''0.0225'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.25%  '
$ws.Range("D50").Value = '''1.920.45'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.49%  '
$ws.Range("D51").Value = '''17.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.48%  '
